$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Not Listed")

# --- Sheet1: update Date_Created / Date_Expired values ---
$ws1.Range("C2").Value = 43788.43747271357
$ws1.Range("D2").Value = 44154.43747271357
$ws1.Range("C3").Value = 43788.43749121563
$ws1.Range("D3").Value = 44154.43749121563
$ws1.Range("C4").Value = 43788.4375108183
$ws1.Range("D4").Value = 44154.4375108183
$ws1.Range("C5").Value = 43788.43753119933
$ws1.Range("D5").Value = 44154.43753119933
$ws1.Range("C6").Value = 43788.43754984346
$ws1.Range("D6").Value = 44154.43754984346
$ws1.Range("C7").Value = 43788.43756985311
$ws1.Range("D7").Value = 44154.43756985311
$ws1.Range("C8").Value = 43788.43759630928
$ws1.Range("D8").Value = 44154.43759630928
$ws1.Range("C9").Value = 43788.43761429244
$ws1.Range("D9").Value = 44154.43761429244
$ws1.Range("C10").Value = 43788.43763458207
$ws1.Range("D10").Value = 44154.43763458207
$ws1.Range("C11").Value = 43788.43765205468
$ws1.Range("D11").Value = 44154.43765205468

# --- Sheet1: fill URL1 column (G) for rows 3-11 ---
$ws1.Range("G3").Value = 'No, individual is not listed'
$ws1.Range("G4").Value = 'No, individual is not listed'
$ws1.Range("G5").Value = 'No, individual is not listed'
$ws1.Range("G6").Value = 'No, individual is not listed'
$ws1.Range("G7").Value = 'No, individual is not listed'
$ws1.Range("G8").Value = 'No, individual is not listed'
$ws1.Range("G9").Value = 'No, individual is not listed'
$ws1.Range("G10").Value = 'No, individual is not listed'
$ws1.Range("G11").Value = 'No, individual is not listed'

# --- Sheet1: update selection ---
$ws1.Range("D13").Select() | Out-Null

# --- Not Listed sheet: populate data ---
$ws3.Range("A2").Value = 'Hello World '
$ws3.Range("B2").Value = 'Hello world '
$ws3.Range("C2").Value = '2019-11-14 15:06:43.151148'
$ws3.Range("D2").Value = '2020-11-14 15:06:43.151148'
$ws3.Range("E2").Value = 'Temple University'
$ws3.Range("F2").Value = 'Phladelphia, PA '
$ws3.Range("A3").Value = 'Achiron'
$ws3.Range("B3").Value = 'Anat'
$ws3.Range("C3").Value = '2019-11-14 15:06:46.292269'
$ws3.Range("D3").Value = '2020-11-14 15:06:46.292269'
$ws3.Range("E3").Value = 'Temple University'
$ws3.Range("F3").Value = 'Phladelphia, PA '
$ws3.Range("A4").Value = 'Afsar'
$ws3.Range("B4").Value = 'Salman'
$ws3.Range("C4").Value = '2019-11-14 15:06:50.515504'
$ws3.Range("D4").Value = '2020-11-14 15:06:50.515504'
$ws3.Range("E4").Value = 'Temple University'
$ws3.Range("F4").Value = 'Phladelphia, PA '
$ws3.Range("A5").Value = 'Akgun'
$ws3.Range("B5").Value = 'Katia'
$ws3.Range("C5").Value = '2019-11-14 15:06:53.411093'
$ws3.Range("D5").Value = '2020-11-14 15:06:53.411093'
$ws3.Range("E5").Value = 'Temple University'
$ws3.Range("F5").Value = 'Phladelphia, PA '
$ws3.Range("A6").Value = 'Alroughani'
$ws3.Range("B6").Value = 'Raed'
$ws3.Range("C6").Value = '2019-11-14 15:06:56.647043'
$ws3.Range("D6").Value = '2020-11-14 15:06:56.647043'
$ws3.Range("E6").Value = 'Temple University'
$ws3.Range("F6").Value = 'Phladelphia, PA '
$ws3.Range("A7").Value = 'Bass'
$ws3.Range("B7").Value = 'Ann'
$ws3.Range("C7").Value = '2019-11-14 15:06:59.691783'
$ws3.Range("D7").Value = '2020-11-14 15:06:59.691783'
$ws3.Range("E7").Value = 'Temple University'
$ws3.Range("F7").Value = 'Phladelphia, PA '
$ws3.Range("A8").Value = 'Berkovich'
$ws3.Range("B8").Value = 'Regina'
$ws3.Range("C8").Value = '2019-11-14 15:07:02.691349'
$ws3.Range("D8").Value = '2020-11-14 15:07:02.691349'
$ws3.Range("E8").Value = 'Temple University'
$ws3.Range("F8").Value = 'Phladelphia, PA '
$ws3.Range("A9").Value = 'Broadley'
$ws3.Range("B9").Value = 'Simon'
$ws3.Range("C9").Value = '2019-11-14 15:07:05.790248'
$ws3.Range("D9").Value = '2020-11-14 15:07:05.790248'
$ws3.Range("E9").Value = 'Temple University'
$ws3.Range("F9").Value = 'Phladelphia, PA '
$ws3.Range("A10").Value = 'Celius'
$ws3.Range("B10").Value = 'Elisabeth'
$ws3.Range("C10").Value = '2019-11-14 15:07:08.776399'
$ws3.Range("D10").Value = '2020-11-14 15:07:08.776399'
$ws3.Range("E10").Value = 'Temple University'
$ws3.Range("F10").Value = 'Phladelphia, PA '
$ws3.Range("A11").Value = 'Hello World '
$ws3.Range("B11").Value = 'Hello world '
$ws3.Range("C11").Value = '2019-11-19 10:29:59.241031'
$ws3.Range("D11").Value = '2020-11-19 10:29:59.241031'
$ws3.Range("E11").Value = 'Temple University'
$ws3.Range("F11").Value = 'Phladelphia, PA '
$ws3.Range("A12").Value = 'Achiron'
$ws3.Range("B12").Value = 'Anat'
$ws3.Range("C12").Value = '2019-11-19 10:30:00.934701'
$ws3.Range("D12").Value = '2020-11-19 10:30:00.934701'
$ws3.Range("E12").Value = 'Temple University'
$ws3.Range("F12").Value = 'Phladelphia, PA '
$ws3.Range("A13").Value = 'Afsar'
$ws3.Range("B13").Value = 'Salman'
$ws3.Range("C13").Value = '2019-11-19 10:30:02.695622'
$ws3.Range("D13").Value = '2020-11-19 10:30:02.695622'
$ws3.Range("E13").Value = 'Temple University'
$ws3.Range("F13").Value = 'Phladelphia, PA '
$ws3.Range("A14").Value = 'Akgun'
$ws3.Range("B14").Value = 'Katia'
$ws3.Range("C14").Value = '2019-11-19 10:30:04.306475'
$ws3.Range("D14").Value = '2020-11-19 10:30:04.306475'
$ws3.Range("E14").Value = 'Temple University'
$ws3.Range("F14").Value = 'Phladelphia, PA '
$ws3.Range("A15").Value = 'Alroughani'
$ws3.Range("B15").Value = 'Raed'
$ws3.Range("C15").Value = '2019-11-19 10:30:06.035309'
$ws3.Range("D15").Value = '2020-11-19 10:30:06.035309'
$ws3.Range("E15").Value = 'Temple University'
$ws3.Range("F15").Value = 'Phladelphia, PA '
$ws3.Range("A16").Value = 'Bass'
$ws3.Range("B16").Value = 'Ann'
$ws3.Range("C16").Value = '2019-11-19 10:30:08.321122'
$ws3.Range("D16").Value = '2020-11-19 10:30:08.321122'
$ws3.Range("E16").Value = 'Temple University'
$ws3.Range("F16").Value = 'Phladelphia, PA '
$ws3.Range("A17").Value = 'Berkovich'
$ws3.Range("B17").Value = 'Regina'
$ws3.Range("C17").Value = '2019-11-19 10:30:09.874867'
$ws3.Range("D17").Value = '2020-11-19 10:30:09.874867'
$ws3.Range("E17").Value = 'Temple University'
$ws3.Range("F17").Value = 'Phladelphia, PA '
$ws3.Range("A18").Value = 'Broadley'
$ws3.Range("B18").Value = 'Simon'
$ws3.Range("C18").Value = '2019-11-19 10:30:11.627891'
$ws3.Range("D18").Value = '2020-11-19 10:30:11.627891'
$ws3.Range("E18").Value = 'Temple University'
$ws3.Range("F18").Value = 'Phladelphia, PA '
$ws3.Range("A19").Value = 'Celius'
$ws3.Range("B19").Value = 'Elisabeth'
$ws3.Range("C19").Value = '2019-11-19 10:30:13.137524'
$ws3.Range("D19").Value = '2020-11-19 10:30:13.137524'
$ws3.Range("E19").Value = 'Temple University'
$ws3.Range("F19").Value = 'Phladelphia, PA '

# --- Workbook-level window tab ratio ---
$excel.ActiveWindow.TabRatio = 0.351
